$word.UserName = "Andreas Bayha"
$d = $word.ActiveDocument
$d.TrackRevisions = $true

# Locate the unique anchor text containing the figure width markdown
# attribute: ".png){:width="300" :class="" and change just the leading
# digit of "300" from "3" to "5" (so the figure becomes wider: 500px).
$anchor = $d.Content
$anchor.Find.Execute('.png){:width="300" :class="', $true, $false, $false, $false, $false, $true, 1, $false, '', 0)

if ($anchor.Find.Found) {
    $matchStart = $anchor.Start
    $digitOffset = $matchStart + 14   # offset of the "3" in `300` within the match
    $digitRange = $d.Range($digitOffset, $digitOffset + 1)

    if ($digitRange.Text -eq "3") {
        $digitRange.Delete()
        $insertionPoint = $d.Range($digitOffset, $digitOffset)
        $insertionPoint.InsertAfter("5")
    }
}
